# Applies cryptocurrency price/volume updates per commit:
# "Updated cryptos list on Tue Sep 24 11:15:50 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.418.69"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3
$ws.Range("D3").Value = "2.641.00"
$ws.Range("E3").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'602.25"
$ws.Range("E5").Value = "  +1.98%  "

# Row 6
$ws.Range("D6").Value = "'146.13"
$ws.Range("E6").Value = "  +1.82%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.588"
$ws.Range("E8").Value = "  +0.45%  "

# Row 9
$ws.Range("D9").Value = "'0.108"
$ws.Range("E9").Value = "  +1.38%  "

# Row 10
$ws.Range("D10").Value = "'5.57"
$ws.Range("E10").Value = "  -0.51%  "

# Row 11
$ws.Range("D11").Value = "'0.369"
$ws.Range("E11").Value = "  +4.61%  "

# Row 12
$ws.Range("E12").Value = "  -0.21%  "

# Row 13
$ws.Range("D13").Value = "'27.48"
$ws.Range("E13").Value = "  +0.58%  "

# Row 14
$ws.Range("D14").Value = "3.117.73"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("D15").Value = "63.280.19"
$ws.Range("E15").Value = "  -0.10%  "

# Row 16
$ws.Range("D16").Value = "'0.0000146"
$ws.Range("E16").Value = "  +0.96%  "

# Row 17
$ws.Range("D17").Value = "2.643.15"
$ws.Range("E17").Value = "  +0.14%  "

# Row 18
$ws.Range("D18").Value = "'11.45"
$ws.Range("E18").Value = "  +0.76%  "

# Row 19
$ws.Range("D19").Value = "'4.56"
$ws.Range("E19").Value = "  +4.74%  "

# Row 20
$ws.Range("D20").Value = "'342.75"
$ws.Range("E20").Value = "  +1.12%  "

# Row 21
$ws.Range("D21").Value = "'6.94"
$ws.Range("E21").Value = "  +3.39%  "

# Row 22
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "  -3.48%  "

# Row 24
$ws.Range("D24").Value = "'66.62"
$ws.Range("E24").Value = "  -0.70%  "

# Row 25
$ws.Range("D25").Value = "'1.69"
$ws.Range("E25").Value = "  +2.34%  "

# Row 26
$ws.Range("D26").Value = "'9.05"
$ws.Range("E26").Value = "  +8.20%  "

# Row 27
$ws.Range("D27").Value = "'575.96"
$ws.Range("E27").Value = "  +6.26%  "

# Row 28
$ws.Range("E28").Value = "  +2.68%  "

# Row 29
$ws.Range("D29").Value = "'0.163"
$ws.Range("E29").Value = "  -1.74%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'7.92"
$ws.Range("E31").Value = "  +2.16%  "

# Row 32
$ws.Range("D32").Value = "'2.05"
$ws.Range("E32").Value = "  +4.79%  "

# Row 33
$ws.Range("E33").Value = "  -3.22%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0820"
$ws.Range("E34").Value = "  +2.03%  "

# Row 35
$ws.Range("D35").Value = "'5.21"
$ws.Range("E35").Value = "  +7.68%  "

# Row 36
$ws.Range("D36").Value = "'167.34"
$ws.Range("E36").Value = "  -4.52%  "

# Row 37
$ws.Range("D37").Value = "'0.406"
$ws.Range("E37").Value = "  +0.97%  "

# Row 38
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.17%  "

# Row 39
$ws.Range("D39").Value = "'1.92"
$ws.Range("E39").Value = "  +7.13%  "

# Row 40
$ws.Range("D40").Value = "'19.10"
$ws.Range("E40").Value = "  +0.51%  "

# Row 41
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("D42").Value = "'168.57"
$ws.Range("E42").Value = "  -1.29%  "

# Row 43
$ws.Range("D43").Value = "'3.76"
$ws.Range("E43").Value = "  +1.01%  "

# Row 44
$ws.Range("D44").Value = "'22.15"
$ws.Range("E44").Value = "  -0.61%  "

# Row 45
$ws.Range("D45").Value = "'0.0569"
$ws.Range("E45").Value = "  +1.15%  "

# Row 46
$ws.Range("D46").Value = "'0.629"
$ws.Range("E46").Value = "  -0.23%  "

# Row 47
$ws.Range("D47").Value = "'0.0245"
$ws.Range("E47").Value = "  +3.02%  "

# Row 48
$ws.Range("D48").Value = "'0.0961"
$ws.Range("E48").Value = "  +0.16%  "

# Row 49
$ws.Range("D49").Value = "'1.88"
$ws.Range("E49").Value = "  +10.93%  "

# Row 50
$ws.Range("D50").Value = "'18.72"
$ws.Range("E50").Value = "  +0.16%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.28"
$ws.Range("E51").Value = "  -0.65%  "
